# This workbook holds a daily-updated "Kiwi" price table.
# A new observation is inserted as the (new) row 110 - the table is kept
# sorted, so every record that used to live at row 110 onward is pushed
# down by one row, and a brand-new record is written into row 110.
#
# Net effect (matching the target diff):
#   - rows 1-109            : unchanged
#   - row 110                : new record (date 44510, etc.)
#   - rows 111-151           : old rows 110-150, shifted down by one
#   - dimension grows from A1:T150 to A1:T151

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new blank row at position 111. This shifts the old rows
# 111..150 down to 112..151, while row 110 (soon to be overwritten with
# new data) stays put for now.
$ws.Rows.Item(111).Insert()

# The newly inserted row 111 is blank; populate it with the data that
# used to sit in row 110 before the edit (the record simply moved down
# one slot).
$ws.Cells.Item(111,1).Value  = 7
$ws.Cells.Item(111,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(111,3).Value  = "Ñuble"
$ws.Cells.Item(111,4).Value  = 44376
$ws.Cells.Item(111,5).Value  = 16
$ws.Cells.Item(111,6).Value  = "Fruta"
$ws.Cells.Item(111,7).Value  = 100101
$ws.Cells.Item(111,8).Value  = "Berries"
$ws.Cells.Item(111,9).Value  = 100101007
$ws.Cells.Item(111,10).Value = "Kiwi"
$ws.Cells.Item(111,11).Value = "Hayward"
$ws.Cells.Item(111,12).Value = "Primera"
$ws.Cells.Item(111,13).Value = 120
$ws.Cells.Item(111,14).Value = 10000
$ws.Cells.Item(111,15).Value = 11000
$ws.Cells.Item(111,16).Value = 10500
$ws.Cells.Item(111,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(111,18).Value = "Provincia de Curicó"
$ws.Cells.Item(111,19).Value = 583
$ws.Cells.Item(111,20).Value = 18

# Now overwrite row 110 with the new observation's values. Columns
# A,B,C,E-L,Q,R,T keep the same value they already had, only the date
# (D), volume (M), min/max/avg price (N/O/P) and $/Kg (S) change.
$ws.Cells.Item(110,4).Value  = 44510
$ws.Cells.Item(110,13).Value = 60
$ws.Cells.Item(110,14).Value = 18000
$ws.Cells.Item(110,15).Value = 19000
$ws.Cells.Item(110,16).Value = 18500
$ws.Cells.Item(110,19).Value = 1028
